$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 35: add Category (E35) ---
$ws.Cells.Item(35, 5).Value2 = "All"

# --- Row 36: add Data (D36) and Category (E36); set row height ---
$ws.Rows.Item(36).RowHeight = 139.5

$d36 = @'
스티어링휠 잠금장치는 차량이 도난된 경우에 차량 조향을 어렵게 만듭니다. 스티어링휠 잠금장치가 잠기거나 해제될 때에는 기계음이 들릴 수 있습니다.
<br><h3>스티어링휠 잠금장치 작동</h3>
<br>차량이 외부에서 잠기고 엔진이 꺼지면 스티어링휠 잠금장치가 작동됩니다. 차량을 잠그지 않은 상태로 두면 잠시 후에 스티어링휠 잠금장치가 자동으로 작동합니다.</br>
<br>
<br><h3>스티어링휠 잠금장치 작동 해제</h3>
<br>차량을 외부에서 잠금 해제하면 스티어링휠 잠금장치가 작동 해제됩니다. 차량이 잠기지 않은 경우에 스티어링휠 잠금장치는 리모컨이 실내에 있고 차량의 시동이 걸린 상태이면 작동하지 않습니다.
'@
$ws.Cells.Item(36, 4).Value2 = $d36
$ws.Cells.Item(36, 4).WrapText = $true
$ws.Cells.Item(36, 5).Value2 = "All"

# --- Row 37: brand new row ---
$ws.Rows.Item(37).RowHeight = 174.75

$a37 = @'
보증 불만으로 인입되는 경우
'@
$ws.Cells.Item(37, 1).Value2 = $a37
$ws.Cells.Item(37, 2).Value2 = 45454
$ws.Cells.Item(37, 3).Value2 = $ws.Cells.Item(36, 3).Value2

$d37 = @'
고객이 서비스센터에 내방했으나, 여러 이유로 보증 적용이 되지 않아 문의를 줄 수 있습니다.
<br>
<br><em><string>서비스센터 현장의 PST의 의견에 따라 진행된 경우 이 의견을 무시하지 마십시오.
<br>또한 이 내용을 VCK로 에스컬레이션 하지 마십시오.</string></em>
<br>
<br>서비스센터로 에스컬레이션하여 PST의 판단으로서 진행된 상황인지 확인합니다.
<br>PST 판단하에 진행된 건인 경우, 고객에게 다음과 같은 안내할 것을 권장합니다.
<br>
<br><em>"보증에 대한 확인을 진행하였으나, 동일하게 답변이 되는 점에 대해 사과드립니다.
<br>해당 내용으로 추가적으로 확인했으나 보증 적용은 어렵습니다.</em>
'@
$ws.Cells.Item(37, 4).Value2 = $d37
$ws.Cells.Item(37, 4).WrapText = $true
$ws.Cells.Item(37, 5).Value2 = "All"

# --- View state: mirror the author's final scroll/selection position ---
$ws.Range("D40").Select()
